# Update the ITI "Backlog" sheet: mark the 21 open items as "Resolvido"
# (was "Pendente"), and move the active selection to M17:M18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# Column I ("Status") for rows 2-22 goes from "Pendente" to "Resolvido".
$ws.Range("I2:I22").Value = "Resolvido"

# Match the author's final selection on the ITI sheet.
$ws.Activate()
$ws.Range("M17:M18").Select()
